$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2799.2666
$ws.Range("J17").Value = 3249.3333
$ws.Range("L17").Value = 9747.999899999999
$ws.Range("N17").Value = -10083.9999
$ws.Range("H116").Value = 9055
$ws.Range("I116").Value = 8949.333000000001
$ws.Range("J116").Value = 9266.333000000001
$ws.Range("K116").Value = 8949.333000000001
$ws.Range("L116").Value = 9266.333000000001
$ws.Range("M116").Value = -5507.333000000001
$ws.Range("N116").Value = -16150.333
$ws.Range("H131").Value = 3745.8125
$ws.Range("I131").Value = 918.8
$ws.Range("K131").Value = 2756.4
$ws.Range("M131").Value = 2283.6
$ws.Range("H132").Value = 1469.4584
$ws.Range("I132").Value = 1488.0526
$ws.Range("K132").Value = 4464.1578
$ws.Range("M132").Value = -1934.1578
$ws.Range("H137").Value = 3365.3076
$ws.Range("J137").Value = 4324.375
$ws.Range("L137").Value = 12973.125
$ws.Range("N137").Value = -18073.125
$ws.Range("H138").Value = 3901.9155
$ws.Range("J138").Value = 4276.9355
$ws.Range("L138").Value = 12830.8065
$ws.Range("N138").Value = -23110.8065
$ws.Range("H141").Value = 3501.8572
$ws.Range("I141").Value = 3316.0908
$ws.Range("J141").Value = 4183
$ws.Range("K141").Value = 9948.2724
$ws.Range("L141").Value = 12549
$ws.Range("M141").Value = -4768.2724
$ws.Range("N141").Value = -22909

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 814.2143
$ws.Range("I2").Value = 708.25
$ws.Range("K2").Value = 708.25
$ws.Range("M2").Value = -595.25
$ws.Range("H32").Value = 18089.44
$ws.Range("I32").Value = 8748.913
$ws.Range("K32").Value = 8748.913
$ws.Range("M32").Value = -8461.913
$ws.Range("H45").Value = 3992.125
$ws.Range("I45").Value = 1919.8572
$ws.Range("K45").Value = 1919.8572
$ws.Range("M45").Value = -1542.8572
$ws.Range("H61").Value = 1404.4546
$ws.Range("I61").Value = 1401.7778
$ws.Range("K61").Value = 1401.7778
$ws.Range("M61").Value = -1189.7778
$ws.Range("H63").Value = 7135
$ws.Range("I63").Value = 4747.5
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 4747.5
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -4061.5
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 7135
$ws.Range("I66").Value = 4747.5
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 23737.5
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -20305.5
$ws.Range("N66").Value = -56864
$ws.Range("H102").Value = 1348.4
$ws.Range("I102").Value = 939.63635
$ws.Range("J102").Value = 2472.5
$ws.Range("K102").Value = 939.63635
$ws.Range("L102").Value = 2472.5
$ws.Range("M102").Value = 682.36365
$ws.Range("N102").Value = -5716.5
$ws.Range("H116").Value = 814.2143
$ws.Range("I116").Value = 708.25
$ws.Range("K116").Value = 708.25
$ws.Range("M116").Value = 1585.75
$ws.Range("H121").Value = 29999
$ws.Range("J121").Value = 29999
$ws.Range("L121").Value = 29999
$ws.Range("N121").Value = -33493
$ws.Range("H136").Value = 1404.4546
$ws.Range("I136").Value = 1401.7778
$ws.Range("K136").Value = 4205.3334
$ws.Range("M136").Value = -1655.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 814.2143
$ws.Range("I3").Value = 708.25
$ws.Range("K3").Value = 708.25
$ws.Range("M3").Value = -594.25
$ws.Range("H94").Value = 839.1111
$ws.Range("I94").Value = 608.6667
$ws.Range("K94").Value = 608.6667
$ws.Range("M94").Value = -157.6667
$ws.Range("H105").Value = 4309
$ws.Range("J105").Value = 6237.25
$ws.Range("L105").Value = 6237.25
$ws.Range("N105").Value = -9731.25
$ws.Range("H134").Value = 3224.818
$ws.Range("I134").Value = 3224.818
$ws.Range("K134").Value = 9674.454000000002
$ws.Range("M134").Value = -7139.454000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 545.7143
$ws.Range("I16").Value = 820
$ws.Range("K16").Value = 820
$ws.Range("M16").Value = -533
$ws.Range("H56").Value = 93
$ws.Range("I56").Value = 93
$ws.Range("K56").Value = 93
$ws.Range("M56").Value = 752
$ws.Range("H76").Value = 7999
$ws.Range("I76").Value = 7999
$ws.Range("K76").Value = 7999
$ws.Range("M76").Value = -7684
$ws.Range("H79").Value = 7999
$ws.Range("I79").Value = 7999
$ws.Range("K79").Value = 7999
$ws.Range("M79").Value = -6907
$ws.Range("H113").Value = 545.7143
$ws.Range("I113").Value = 820
$ws.Range("K113").Value = 820
$ws.Range("M113").Value = 1350
$ws.Range("H122").Value = 7492.636
$ws.Range("I122").Value = 7166.1333
$ws.Range("K122").Value = 21498.3999
$ws.Range("M122").Value = -19048.3999
$ws.Range("H134").Value = 4742.25
$ws.Range("I134").Value = 3499.5
$ws.Range("J134").Value = 5156.5
$ws.Range("K134").Value = 10498.5
$ws.Range("L134").Value = 15469.5
$ws.Range("M134").Value = -7963.5
$ws.Range("N134").Value = -20539.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111184.78
$ws.Range("I2").Value = 200045.2
$ws.Range("J2").Value = 109.25
$ws.Range("K2").Value = 1200271.2
$ws.Range("L2").Value = 655.5
$ws.Range("M2").Value = -1200158.2
$ws.Range("N2").Value = -881.5
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -2827
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 5750
$ws.Range("J22").Value = 5750
$ws.Range("L22").Value = 17250
$ws.Range("N22").Value = -17588
$ws.Range("H27").Value = 5750
$ws.Range("J27").Value = 5750
$ws.Range("L27").Value = 17250
$ws.Range("N27").Value = -17454
$ws.Range("H69").Value = 4337.3335
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622
$ws.Range("H72").Value = 4337.3335
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112
$ws.Range("H131").Value = 4134.25
$ws.Range("J131").Value = 5157.4
$ws.Range("L131").Value = 15472.2
$ws.Range("N131").Value = -25552.2
$ws.Range("H136").Value = 19286.666
$ws.Range("I136").Value = 11526.667
$ws.Range("K136").Value = 34580.001
$ws.Range("M136").Value = -29480.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14795
$ws.Range("I57").Value = 4659.3335
$ws.Range("K57").Value = 4659.3335
$ws.Range("M57").Value = -3839.3335
$ws.Range("H107").Value = 890.5417
$ws.Range("I107").Value = 656.6
$ws.Range("K107").Value = 656.6
$ws.Range("M107").Value = 1263.4
$ws.Range("H132").Value = 2991.7097
$ws.Range("I132").Value = 2638.64
$ws.Range("K132").Value = 7915.92
$ws.Range("M132").Value = -5385.92
$ws.Range("H136").Value = 25703.182
$ws.Range("J136").Value = 25703.182
$ws.Range("L136").Value = 77109.546
$ws.Range("N136").Value = -82209.546

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2099.5
$ws.Range("I40").Value = 2099.5
$ws.Range("K40").Value = 2099.5
$ws.Range("M40").Value = -1963.5
$ws.Range("H93").Value = 264.85715
$ws.Range("I93").Value = 225.66667
$ws.Range("K93").Value = 225.66667
$ws.Range("M93").Value = 1022.33333
$ws.Range("H122").Value = 6000.75
$ws.Range("I122").Value = 6000.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18002.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15552.25
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4465.911
$ws.Range("I132").Value = 3951.0312
$ws.Range("K132").Value = 11853.0936
$ws.Range("M132").Value = -9323.0936
$ws.Range("H136").Value = 3047.04
$ws.Range("I136").Value = 2917.1365
$ws.Range("K136").Value = 8751.4095
$ws.Range("M136").Value = -6201.4095

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2004000
$ws.Range("J5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("N5").Value = -8224
$ws.Range("H46").Value = 416292.34
$ws.Range("J46").Value = 416292.34
$ws.Range("L46").Value = 416292.34
$ws.Range("N46").Value = -416754.34
$ws.Range("H81").Value = 1423.8667
$ws.Range("I81").Value = 1399.1
$ws.Range("J81").Value = 1473.4
$ws.Range("K81").Value = 2798.2
$ws.Range("L81").Value = 2946.8
$ws.Range("M81").Value = -1737.2
$ws.Range("N81").Value = -5068.8
$ws.Range("H84").Value = 1423.8667
$ws.Range("I84").Value = 1399.1
$ws.Range("J84").Value = 1473.4
$ws.Range("K84").Value = 13991
$ws.Range("L84").Value = 14734
$ws.Range("M84").Value = -8687
$ws.Range("N84").Value = -25342
$ws.Range("H134").Value = 416292.34
$ws.Range("J134").Value = 416292.34
$ws.Range("L134").Value = 1248877.02
$ws.Range("N134").Value = -1253947.02
$ws.Range("H136").Value = 54414.473
$ws.Range("I136").Value = 1208.7858
$ws.Range("K136").Value = 3626.3574
$ws.Range("M136").Value = -1076.3574
